# Append top-50 high-retweet / low-sentiment rows (17-51) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(202, 'Really feeling for people on #newstart tonight, regardless of the Labor hedge position, the Coalition says it’s flat out opposed. #ausvotes', -0.710941856727004, 'Brunswick, Melbourne'),
    @(193, 'Is it just me, but the more I hear, read and see, the more we look like a corrupt third world country. What’s far worse is how many people are being fooled by believing this crook mob, and are going to vote for them to get back into government. Oh Australia! #auspol', -0.9028009162284436, 'Sydney, New South Wales'),
    @(171, '❓How can "Liberals" justify giving Clive Palmer an armchair ride into the Senate when Palmer is:
🏝️In Fiji during the election
👎Guilty of calling Chinese people "mongrels" and "bastards"
🤪Running freakshow candidates with faked credentials
#ausvotes #auspol https://t.co/7c3yJfxm8m', -0.8039950579404832, 'Authorised: Kimberley Kitching, Australian Labor Party, Melbourne, Victoria'),
    @(164, 'Gday Paul. Terrorist Watchlist?
Your honour, I please guilty to assault with a deadly metaphor. 
#ausvotes #auspol
Paul must work for Intelligence. https://t.co/qiRjHurRr6', -0.7764248400926591, 'Sydney, New South Wales'),
    @(146, 'In 2013 the narrative was that Labor had to go because of 6 years of chaos, 2 PMs and racking up debt. In the next 6 years, the Liberals have had even more chaos, even more PMs and even more debt but were rewarded for it. What gives? #AusVotes', -0.8072726605460049, 'Rockhampton, Queensland'),
    @(141, 'If you think #Barnaby Christensen &amp; Angus Taylor were crook before this...you watch em go now. They just got the green light they wanted and will loot us all for everything. #auspol', -0.8524082070216538, 'The Adelaide Diaspora'),
    @(137, '#abc730 Morrison is driving me 😩 insane  he shouts too much.. He''s rude.. He''s a bully.. Ffs who could stand this for another 3 years
People do your thing #auspol #AusVotes2019', -0.924999024719, 'brisbane'),
    @(129, 'I wonder why Chloe Shorten didn’t want to speak to Melissa Fyfe? 🤔 Could it be because Fyfe was obviously intent on writing the most snide and nasty profile possible, despite not being able to find one single source to speak a bad word about a good woman? #AusVotes2019 #auspol https://t.co/db8xZWTbj8', -0.797601329162717, 'Melbourne'),
    @(125, 'Leigh Sales: "If Newspoll has it so wrong, isn''t it time you stop sacking your leaders over it?" #ausvotes', -0.7683270666748285, 'Sydney, Australia'),
    @(123, 'Noticed today some people tweeting they vote  ALP because  #BobHawke has passed  on...sad but Hawke  can''t help u and that vote....although a nice gesture is a vote for someone who  can''t  help u... Bill Shorten will screw u...your kids..your way of life.. the country #ausvotes', -0.8442283431068064, 'Geelong, Victoria'),
    @(119, 'Voting for the major parties in Australia''s federal election tomorrow is the equivalent of picking between death by a thousand cuts and being pecked to death by a chicken.
One might kill the nation slower than the other, but we''re going down either way.', -0.7327415943145753, 'Melbourne, Victoria'),
    @(114, 'Bigotry &amp; racism emboldened. 😡#auspol https://t.co/W640NBP829', -0.8997668339870869, 'Hunter Valley NSW'),
    @(108, '#auspol #ausvotes #ausvotes2019
What is it with the haters? There must be a statistical parallel with anger/hate and a toxic inability to use punctuation and correct spelling.
Is it due to an empty fear driven by fictional threats of what a non-Liberal world would look like?', -0.8174341507256031, 'Sydney, New South Wales'),
    @(107, 'Fuck!  #auspol #AusVotes2019 
Former Australian prime minister Bob Hawke dies aged 89 – rolling coverage
https://t.co/cHaey0CBrR', -0.7940215962007642, 'brisbane'),
    @(107, 'This is about as low as political debate can get: Scott Morrison’s personal beliefs are being used against him as Labor seeks to weaponise faith, as well as mobilising the worst aspects of the SSM debate, writes Dennis Shanahan #auspol #ausvotes https://t.co/8787c3R91x', -0.8843950503505766, 'Sydney, Australia'),
    @(105, '‘Your beliefs are a joke’: footage of Greens candidate Connor Parissis shouting down Christians emerges #auspol #ausvotes https://t.co/qjYcgPBCJK', -0.8616530066356064, 'Sydney, Australia'),
    @(99, 'journalists at the Australian are uncomfortable with the paper’s cheerleading for the Coalition. Wow it''s worse at the Courier-Mail I''ve been told from someone, it''s almost at boiling point. Almost ready for a walkout.  Not everyone who works there votes liberal  ffs #auspol', -0.90053202630952, 'brisbane'),
    @(94, 'I''m always perplexed by these politicians, like Simon Birmingham and Scott Morrison, who insist that governments can''t solve problems. Bit of a profound self-own, isn''t it?
#qanda #ausvotes', -0.8658048901706934, 'Melbourne, Victoria'),
    @(93, 'If Labor does stay ahead and win election, it''ll be most massive clarsterfark in a dumpster fire #Auspol has ever seen. Shorten is so clearly a puppet of unions, globalists. He''s the hollowest of hollowmen. Spineless dolt''d make even Gillard look visionary and principled!', -0.8517128601670266, 'Sydney, New South Wales'),
    @(90, 'Abbott doesn''t like Get up campaigning against him,  says it''s feeling pretty personal. This is from the man who set up a slush fund to get rid of Pauline Hanson #auspol', -0.7316501131281258, 'Ipswich, Qld'),
    @(90, 'I was so excited to vote in my first federal election this year, but now I just feel helpless and totally unheard. As a young person, the election results tonight are highly disappointing and make me worried for the future of the planet 💔#ausvotes #Election2019Results #politas', -0.8570144344121218, 'Tasmania, Australia'),
    @(90, 'The Murdoch rags go from bad to worse lets start with the racist pun and complete it with the Premier of QLD lined up in the crosshairs of a gun.
Totally irresponsible
#Auspol #AdShame 👇👇🏿👇🏻👇🏾 via @Bloss55 https://t.co/6okluX15WC', -0.9206178151071072, 'Australia'),
    @(88, 'Reasons not to vote for Scott Morrison: 
-Voted against Banking Royal Commission 26 times
-Abstained from the same sex marriage vote 
-Voted to cut penalty rates 
- Urged the shadow cabinet to capitalise on “growing concerns of Muslims”
I could go all day...
#auspol #ausvotes', -0.7185311745852233, 'Canberra, Australian Capital Territory'),
    @(88, 'Well that is the end of honest election campaigns #auspol #auspol2019 https://t.co/yk6HAQsROl', -0.8009242350235583, 'Sunshine Coast, Queensland'),
    @(87, 'Turns out, ⁦@JoshFrydenberg⁩''s budget didn''t even last a month. 
#auspol #ausvotes
Frydenberg’s budget forecasts downgraded, RBA flags rate cut https://t.co/XMjEaq3t6j', -0.7705959854647517, 'Brisbane, Australia'),
    @(84, 'This is 2019... I thought Aus was better than this... Why do our people not care about equality, climate change, and kindness?! Why is Aus supporting racism, white supremacy, unequal opportunities, and hatred? RIP Australia 💔 #AustraliaDecides #AUSVote2019 #Election2019Results', -0.9472382236272096, 'Australia'),
    @(84, '“Voters in the marginal seat of Corangamite have been the target of perhaps the biggest pork-barrelling campaign in modern electoral history, showered with the equivalent of $26,500 each if they vote Liberal.” @swrighteconomy @Jess_Irvine #auspol https://t.co/dsucwk3xRr', -0.7055868562310935, 'Australian Capital Territory'),
    @(78, 'Coalition’s plan for public schools over the next decade is $14 billion in cuts, ensuring almost every single public school will fall short of the schooling resource standard. That’s the real outcome we can’t afford. #ChangeTheGovernment #ChangeTheRules #FairFundingNow #ausvotes https://t.co/9NCxszGlYl', -0.8578356597572565, 'Adelaide, South Australia'),
    @(77, 'If Simon Birmingham and the Liberal-Nationals think that handing $500 million to a bunch of rich spivs has "saved the Great Barrier Reef", this country is not going to see the next century. #qanda #ausvotes', -0.8186883935704827, 'Melbourne, Victoria'),
    @(76, 'The repressive new dark age started when Ruper Murdoch &amp; Roger Ailes decided that pitting people against one another was profitable and pursued it with greedy passion.
Religion is just another point of difference to exploit.
God SAVE Australia from Murdoch.
#Auspol https://t.co/8zUYAs2gYH', -0.7629188392311336, 'Australia'),
    @(76, 'It’d be a terrible outcome if the focus on a likely Labor victory in the House of Representatives deflected attention from the very real risk of a Senate controlled by #RWNJ’s. @greens holding the balance of power in the Senate would be a much better outcome #auspol #Election2019', -0.715309275314212, 'Victoria Australia'),
    @(69, 'While everyone is screaming about a coal mine that hasn’t been able to start yet for SIX YEARS, the Liberals are cutting funding to ABC, SBS, SafeWork, NDIS and other departments. Which means STAFF CUTS. #Auspol #AusVotes2019 https://t.co/g1z6howXZs', -0.8399765212088823, 'Rockhampton'),
    @(67, 'Clive Palmer''s $60m election ad blitz appears to have failed to win a single seat https://t.co/zfKctJJ2ff', -0.8083290942013265, 'Sydney/Canberra'),
    @(66, 'Scott Morrison''s claim Coalition saved Great Barrier Reef condemned as ''ridiculous'' https://t.co/Gi4hWaqIjp #AusVotes2019 #ausvotes', -0.71904965210706, 'Sydney/Canberra'),
    @(65, 'The Coalition has announced a further $1.5bn in cuts to the public service to pay for its $1.4bn in election promises, a move that could reduce jobs by 3,000.
#auspol
#ausvotes19
https://t.co/GW7o3OES1p', -0.7495171139016747, 'Sydney')
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}

